$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 data: happy flow for RAYA -> NCR, date range split into days
$ws.Range("B2").Value = "00000690"
$ws.Range("C2").Value = [DateTime]"2020-09-11"
$ws.Range("D2").Value = [DateTime]"2020-11-11"
$ws.Range("E2").Value = "NCR"
$ws.Range("F2").Value = "Success"

# Update the active selection to reflect the new working cell
$ws.Range("E8").Select()

# Set page orientation for the worksheet
$ws.PageSetup.Orientation = 1
